# Updated cryptos list on Sun Mar  3 03:49:23 UTC 2024 with GitHub Actions
#
# The "Price" column (D) stores plain numeric-looking text (e.g. "410.92",
# "1.00", "0.0500") as text, not numbers, in the source workbook (every
# cell in the sheet is t="inlineStr"). Excel's COM layer auto-converts a
# plain numeric string assigned to .Value into a real number (and even
# silently rounds/renders trailing zeros differently, e.g. "1.00" -> 1),
# which would corrupt the intended text formatting. To avoid that, each
# Price cell being rewritten with a number-like string is temporarily
# forced to Text number format ("@") before the new value is written, and
# the formatting is cleared back to the sheet's default immediately after
# so no stray cell-level styling is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "62.034.19"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.430.57"
$ws.Range("E3").Value = "  -0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "410.92"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6 - Solana
Set-TextValue "D6" "130.27"
$ws.Range("E6").Value = "  +0.26%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.63%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.66%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.28%  "

# Row 11 - Avalanche
Set-TextValue "D11" "43.67"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  +15.79%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +4.95%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.972.38"
$ws.Range("E14").Value = "  -0.40%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.12%  "

# Row 16 - Chainlink
$ws.Range("E16").Value = "  +3.44%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.427.51"
$ws.Range("E17").Value = "  -0.48%  "

# Row 18 - Uniswap
Set-TextValue "D18" "12.37"
$ws.Range("E18").Value = "  +4.88%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +2.73%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "61.963.65"
$ws.Range("E20").Value = "  -0.54%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "519.85"
$ws.Range("E21").Value = "  +28.19%  "

# Row 22 - Litecoin
Set-TextValue "D22" "93.06"
$ws.Range("E22").Value = "  +4.07%  "

# Row 23 - ImmutableX
$ws.Range("E23").Value = "  +3.85%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  +0.50%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "3.35"
$ws.Range("E25").Value = "  +3.68%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "35.03"
$ws.Range("E26").Value = "  +9.03%  "

# Row 27 - Filecoin
Set-TextValue "D27" "9.29"
$ws.Range("E27").Value = "  +8.79%  "

# Row 28 - RenderToken
Set-TextValue "D28" "7.65"
$ws.Range("E28").Value = "  -0.79%  "

# Row 29 - Cosmos
Set-TextValue "D29" "12.16"
$ws.Range("E29").Value = "  +2.72%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -1.74%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -1.82%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -2.69%  "

# Row 33 - InjectiveProtocol
Set-TextValue "D33" "41.95"
$ws.Range("E33").Value = "  -5.02%  "

# Row 34 - was OKB, now Dai (rows 34/35 swapped order)
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D34" "1.00"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35 - was Dai, now OKB
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D35" "58.68"
$ws.Range("E35").Value = "  +11.93%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0500"
$ws.Range("E36").Value = "  +1.22%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 - Stellar
Set-TextValue "D38" "0.139"
$ws.Range("E38").Value = "  +5.49%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +2.18%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +18.56%  "

# Row 41 - Monero
Set-TextValue "D41" "148.05"
$ws.Range("E41").Value = "  +5.22%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +0.71%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  +6.80%  "

# Row 44 - TheGraph
Set-TextValue "D44" "0.318"
$ws.Range("E44").Value = "  +1.47%  "

# Row 45 - NEARProtocol
$ws.Range("E45").Value = "  +7.13%  "

# Row 46 - ThetaToken
$ws.Range("E46").Value = "  +21.99%  "

# Row 47 - Celestia
Set-TextValue "D47" "16.69"
$ws.Range("E47").Value = "  -0.65%  "

# Row 48 - was EnergySwap, now BitcoinSV (rows 48/49 swapped order)
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D48" "122.46"
$ws.Range("E48").Value = "  +26.40%  "

# Row 49 - was BitcoinSV, now EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "23.18"
$ws.Range("E49").Value = "  +4.71%  "

# Row 50 - Cronos
Set-TextValue "D50" "0.146"
$ws.Range("E50").Value = "  +18.88%  "

# Row 51 - Maker
Set-TextValue "D51" "2.144.83"
$ws.Range("E51").Value = "  +0.97%  "
